# Multiple Users Template.xlsx edit:
#  - The "Email" header column moved from the first column to the last
#    column, and was renamed "Email (Optional)" (Email is now optional
#    because a new "Organizations" drop-down validation column is used
#    instead as the primary identifier when present).
#  - All other header columns shifted one position to the left.
#  - Column A was widened to fit the new, longer "Email (Optional)" text.
#  - The active cell selection moved from E11 to C8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-write the header row (row 1) with the new column order.
$ws.Range("A1").Value = "Email (Optional)"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Organizations (Optional)"
$ws.Range("D1").Value = "Address (Optional)"
$ws.Range("E1").Value = "Gender (Optional)"
$ws.Range("F1").Value = "Birthdate (Optional)"
$ws.Range("G1").Value = "SchoolUniversityJob (Optional)"
$ws.Range("H1").Value = "GpsLocation (Optional)"
$ws.Range("I1").Value = "NationalID (Optional)"
$ws.Range("J1").Value = "MentorName (Optional)"
$ws.Range("K1").Value = "FirstMobile (Optional)"
$ws.Range("L1").Value = "SecondMobile (Optional)"
$ws.Range("M1").Value = "FatherMobile (Optional)"
$ws.Range("N1").Value = "MotherMobile (Optional)"

# Column A now holds the longer "Email (Optional)" text, so widen it to fit.
$ws.Columns.Item(1).ColumnWidth = 17.43

# Move the active selection to C8.
$ws.Range("C8").Select()
